$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2-17
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05)
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
